$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.583.42"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.786.71"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.07"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.77"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +7.26%  "
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("E10").Value = "  +1.67%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0937"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.60%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.043.68"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.70%  "
$ws.Range("E13").Value = "  +11.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.789.22"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +0.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.582.28"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.49%  "
$ws.Range("E16").Value = "  +0.60%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.29"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.39"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "252.82"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0770"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +4.17%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.03%  "
$ws.Range("E22").Value = "  +1.28%  "
$ws.Range("E23").Value = "  +0.80%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.14"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "158.37"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.30"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("E27").Value = "  +2.15%  "
$ws.Range("E28").Value = "  -0.20%  "
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("E30").Value = "  -1.84%  "
$ws.Range("E31").Value = "  +0.11%  "
$ws.Range("E32").Value = "  -0.64%  "
$ws.Range("E33").Value = "  +0.32%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.84"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.440.55"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.88%  "
$ws.Range("E36").Value = "  -1.16%  "
$ws.Range("E37").Value = "  +2.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.630"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.79%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "82.90"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.41%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.80"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +3.83%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.35"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.20%  "
$ws.Range("E42").Value = "  +1.74%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.06"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0506"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("E45").Value = "  +2.96%  "
$ws.Range("E46").Value = "  -2.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.939.43"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "104.24"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +6.80%  "
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("E50").Value = "  +1.41%  "
$ws.Range("E51").Value = "  +5.36%  "
